# Added ifoCAST full series evaluation: for each data row (2..16) the
# naive QoQ error series is re-aligned one quarter later. The value that
# used to sit in column B (the earliest quarter) is dropped, every other
# value shifts one column to the left, and - where a new quarter's error
# is now available (rows 2..6) - a freshly computed value is appended in
# the vacated last column. For the remaining rows (7..16) there is no new
# data yet, so the row simply becomes one cell shorter.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New trailing value appended for rows that now have a newly matched
# ifoCAST quarter (column K in each case).
$newTailValues = @{
    2 = 0.08445119490591374
    3 = -0.1050779264540507
    4 = 0.2888923924969555
    5 = -0.1448629740152609
    6 = 1.455535253737389
}

# Number of populated data cells (columns B..) in each row *before* the edit.
$rowLengths = @{
    2 = 10; 3 = 10; 4 = 10; 5 = 10; 6 = 10; 7 = 10
    8 = 9; 9 = 8; 10 = 7; 11 = 6; 12 = 5; 13 = 4; 14 = 3; 15 = 2; 16 = 1
}

foreach ($row in 2..16) {
    $len = $rowLengths[$row]
    # Column B is index 2, so the last populated column index is 1 + len.
    $lastColIndex = 1 + $len

    # Shift every value one column to the left: new col c = old col c+1.
    for ($colIndex = 2; $colIndex -lt $lastColIndex; $colIndex++) {
        $ws.Cells.Item($row, $colIndex).Value = $ws.Cells.Item($row, $colIndex + 1).Value2
    }

    if ($newTailValues.ContainsKey($row)) {
        # A new quarter's error value is now available - place it in the
        # vacated last column.
        $ws.Cells.Item($row, $lastColIndex).Value = $newTailValues[$row]
    } else {
        # No new value yet - the row shrinks by one cell.
        $ws.Cells.Item($row, $lastColIndex).ClearContents()
    }
}
